# Fruta / hortaliza, semanal
# Insert a new weekly record into the Zanahoria (carrot) price table.
# The new observation is inserted at row 179, pushing the existing
# rows 179-210 down to 180-211.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 179 (shifts rows 179:210 down to 180:211,
# row 211 receives what used to be row 210's data).
$ws.Rows("179:179").Insert()

# Populate the newly inserted row 179 with the new weekly record.
# This mirrors the previous row 179's values except for the date
# (column D) and the origin (column O), which are the actual new data.
$ws.Range("A179").Value = 5
$ws.Range("B179").Value = "Macroferia Regional de Talca"
$ws.Range("C179").Value = "Maule"
$ws.Range("D179").Value = 44504
$ws.Range("E179").Value = 7
$ws.Range("F179").Value = 100114013
$ws.Range("G179").Value = "Zanahoria"
$ws.Range("H179").Value = "Sin especificar"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 400
$ws.Range("K179").Value = 8000
$ws.Range("L179").Value = 8000
$ws.Range("M179").Value = 8000
$ws.Range("N179").Value = "`$/saco 20 kilos"
$ws.Range("O179").Value = "Provincia del Elquí"
$ws.Range("P179").Value = 400
$ws.Range("Q179").Value = 20
$ws.Range("R179").Value = "Hortaliza"
